$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("F17").Value = "2026-02-13T17:22:29"
$ws.Range("G17").ClearContents()

# Row 38
$ws.Range("F38").Value = "2026-02-14T10:34:29"

# Row 42
$ws.Range("F42").Value = "2026-02-14T23:58:34"

# Row 87
$ws.Range("F87").Value = "2026-02-14T13:56:54"

# Row 117
$ws.Range("F117").Value = "2026-02-14T12:54:38"

# Row 121
$ws.Range("G121").ClearContents()

# Row 122
$ws.Range("F122").Value = "2026-02-14T10:01:33"

# Row 123
$ws.Range("F123").Value = "2026-02-18T07:10:15"

# Row 124
$ws.Range("F124").Value = "2026-02-14T23:07:40"

# Row 126
$ws.Range("G126").ClearContents()

# Row 127
$ws.Range("F127").Value = "2026-02-14T10:01:33"

# Row 128
$ws.Range("F128").Value = "2026-02-14T10:01:59"

# Row 130
$ws.Range("F130").Value = "2026-02-14T11:49:08"
$ws.Range("G130").ClearContents()

# Row 131
$ws.Range("F131").Value = "2026-02-14T09:11:41"

# Row 132
$ws.Range("F132").Value = "2026-02-18T06:53:30"

# Row 134
$ws.Range("F134").Value = "2026-02-14T11:49:08"

# Row 136
$ws.Range("G136").ClearContents()

# Row 137
$ws.Range("F137").Value = "2026-02-14T09:11:41"
$ws.Range("G137").ClearContents()

# Row 144
$ws.Range("F144").Value = "2026-02-14T09:11:52"

# Row 153
$ws.Range("F153").Value = "2026-02-14T09:59:02"

# Row 155
$ws.Range("F155").Value = "2026-02-14T09:20:30"

# Row 167
$ws.Range("F167").Value = "2026-02-14T20:37:37"

# Row 168
$ws.Range("F168").Value = "2026-02-14T16:23:01"

# Row 172
$ws.Range("F172").Value = "2026-02-02T20:08:57"
$ws.Range("G172").Value = "PA_054"

# Row 173
$ws.Range("F173").Value = "2026-02-11T10:28:00"
$ws.Range("G173").Value = "PA_043"

# Row 184
$ws.Range("F184").Value = "2026-02-14T02:37:16"

# Row 257
$ws.Range("F257").Value = "2026-02-14T16:23:01"

# Row 271
$ws.Range("G271").ClearContents()

# Row 280
$ws.Range("F280").Value = "2026-02-14T17:17:08"

# Row 306
$ws.Range("F306").Value = "2026-02-13T15:55:33"

# Row 310
$ws.Range("F310").Value = "2026-02-14T12:38:01"

# Row 326
$ws.Range("F326").Value = "2026-02-14T11:09:27"
